$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("E3").Value = 16.199
$ws.Range("E4").Value = 16.495
$ws.Range("A11").Value = -21.705
$ws.Range("A12").Value = -21.657
$ws.Range("E14").Value = 16.905
$ws.Range("A15").Value = -22.044
$ws.Range("E26").Value = 16.373
$ws.Range("A27").Value = -21.814
$ws.Range("A28").Value = -21.777
$ws.Range("A31").Value = -21.849
$ws.Range("E31").Value = 16.2
$ws.Range("A32").Value = -21.684
$ws.Range("E35").Value = 16.545
$ws.Range("A36").Value = -20.339
$ws.Range("E37").Value = 16.643
$ws.Range("A38").Value = -19.69600000000001
$ws.Range("E39").Value = 16.376
$ws.Range("E40").Value = 16.547
$ws.Range("E45").Value = 16.851
$ws.Range("A46").Value = -21.802
$ws.Range("E52").Value = 16.817
$ws.Range("A54").Value = -21.764
$ws.Range("A55").Value = -22.214
$ws.Range("A56").Value = -22.001
$ws.Range("E57").Value = 16.643
$ws.Range("A67").Value = -21.565
$ws.Range("A69").Value = -21.721
$ws.Range("A72").Value = -21.445
$ws.Range("A73").Value = -19.981
$ws.Range("E81").Value = 16.709
$ws.Range("A83").Value = -21.905
$ws.Range("E83").Value = 16.554
$ws.Range("A86").Value = -22.093
$ws.Range("A91").Value = -21.652
$ws.Range("A93").Value = -21.259
$ws.Range("A99").Value = -20.54
$ws.Range("E100").Value = 16.568
$ws.Range("E102").Value = 16.52
